{"js": "// Office.js (Word JavaScript API) edit script.\n// Body: `async (context) => { ... }`\n//\n// Target edit: the paragraph's sentence\n//   \"To remove some weird bugs, remove the collisions of the meshes.\"\n// becomes\n//   \"To remove some weird bugs, remove the collisions of the meshes or disable gravity for some meshes.\"\n// i.e. the trailing period is replaced with\n//   \" or disable gravity for some meshes.\"\n\nconst oldText = \"To remove some weird bugs, remove the collisions of the meshes.\";\nconst newText = \"To remove some weird bugs, remove the collisions of the meshes or disable gravity for some meshes.\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet applied = false;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text.indexOf(oldText) !== -1) {\n    if (text === oldText) {\n      // Whole paragraph is exactly the sentence: replace the paragraph's range.\n      para.getRange().insertText(newText, \"Replace\");\n    } else {\n      // Sentence is part of a larger paragraph: replace just the matching text.\n      const found = para.search(oldText, { matchCase: true });\n      found.load(\"items\");\n      await context.sync();\n      for (let j = 0; j < found.items.length; j++) {\n        found.items[j].insertText(newText, \"Replace\");\n      }\n    }\n    applied = true;\n  }\n}\n\nif (!applied) {\n  // Fallback: search the whole document body in case the sentence spans\n  // differently than expected (defensive, should not normally trigger).\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  for (let j = 0; j < found.items.length; j++) {\n    found.items[j].insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Target edit: the paragraph's sentence\n#   \"To remove some weird bugs, remove the collisions of the meshes.\"\n# becomes\n#   \"To remove some weird bugs, remove the collisions of the meshes or disable gravity for some meshes.\"\n# i.e. the trailing period is replaced with\n#   \" or disable gravity for some meshes.\"\n\n$d = $word.ActiveDocument\n\n$oldText = \"To remove some weird bugs, remove the collisions of the meshes.\"\n$newText = \"To remove some weird bugs, remove the collisions of the meshes or disable gravity for some meshes.\"\n\n# Use Find/Replace over the whole document story so the run's character\n# formatting (Arial, 12pt) is preserved automatically.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n\n# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# Wrap=1 -> wdFindContinue, Replace=2 -> wdReplaceAll\n$find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n"}
